{"js": "// The published site rebuild drops the page-footer boilerplate paragraphs\n// (the blank spacer, the \"Ver no Jupiter...\" line, and the \"\u00a9 2020...\"\n// copyright line) that used to follow the last \"Requisitos\" entry.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targets = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\n// Find the \"LOM3229: ...\" requirement paragraph; the three paragraphs to\n// remove are the ones immediately following it (blank spacer, \"Ver no\n// Jupiter...\" line, \"\u00a9 2020...\" line).\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOM3229\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (anchorIndex !== -1) {\n  for (let i = anchorIndex + 1; i < items.length && toDelete.length < targets.length; i++) {\n    const text = items[i].text.replace(/[\\r\\x07]+$/, \"\");\n    if (text === targets[toDelete.length]) {\n      toDelete.push(items[i]);\n    } else {\n      break;\n    }\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# The published site rebuild drops the page-footer boilerplate paragraphs\n# (the blank spacer, the \"Ver no Jupiter...\" line, and the \"\u00a9 2020...\"\n# copyright line) that used to follow the last \"Requisitos\" entry\n# (\"LOM3229: ...\").\n$d = $word.ActiveDocument\n\n$targets = @(\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n# Locate the \"LOM3229: ...\" requirement paragraph.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\")\n    if ($t.IndexOf(\"LOM3229\") -ge 0) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ge 0) {\n    # Confirm the three paragraphs immediately following the anchor are the\n    # expected boilerplate before removing them, then delete starting from\n    # the last one so earlier indices stay valid.\n    $matchCount = 0\n    for ($k = 0; $k -lt $targets.Length; $k++) {\n        $idx = $anchorIndex + 1 + $k\n        if ($idx -gt $d.Paragraphs.Count) { break }\n        $t = $d.Paragraphs.Item($idx).Range.Text.TrimEnd(\"`r\")\n        if ($t -ne $targets[$k]) { break }\n        $matchCount++\n    }\n\n    for ($k = $matchCount - 1; $k -ge 0; $k--) {\n        $idx = $anchorIndex + 1 + $k\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
